# The test case formula used to live in A2; it now needs to live in A3
# (A2's neighbor cells are needed for the new double-border support), so
# move the formula down one row: clear A2 and write the same formula into A3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").ClearContents()
$ws.Range("A3").Formula = "=A1+1"

# Keep the selection in sync with the relocated cell.
$ws.Range("A3").Select()
